$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"0.002529333333333334"
$ws.Range("H2").Value = [double]"0.007588"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.02564166666666666"
$ws.Range("N2").Value = [double]"0.07692499999999999"
$ws.Range("O2").Value = [double]"0.0006780701807970013"
$ws.Range("P2").Value = [double]"0.0006780701807970013"
$ws.Range("Q2").Value = [double]"6.485632222222222E-05"
$ws.Range("R2").Value = [double]"0.0005837069"
$ws.Range("S2").Value = [double]"0.0006780701807970013"
$ws.Range("T2").Value = [double]"0.0006780701807970013"

# Row 3
$ws.Range("G3").Value = [double]"0.002529333333333334"
$ws.Range("H3").Value = [double]"0.007588"
$ws.Range("M3").Value = [double]"0.01112833333333333"
$ws.Range("O3").Value = [double]"0.0002942784918545062"
$ws.Range("P3").Value = [double]"0.0002942784918545062"
$ws.Range("Q3").Value = [double]"2.814726444444444E-05"
$ws.Range("S3").Value = [double]"0.0002942784918545062"
$ws.Range("T3").Value = [double]"0.0002942784918545062"

# Row 4
$ws.Range("G4").Value = [double]"0.002529333333333334"
$ws.Range("H4").Value = [double]"0.007588"
$ws.Range("O4").Value = [double]"0.0009221225577320236"
$ws.Range("P4").Value = [double]"0.0009221225577320235"
$ws.Range("Q4").Value = [double]"8.819953955555556E-05"
$ws.Range("R4").Value = [double]"0.000793795856"
$ws.Range("S4").Value = [double]"0.0009221225577320236"
$ws.Range("T4").Value = [double]"0.0009221225577320235"

# Row 5
$ws.Range("G5").Value = [double]"0.002529333333333334"
$ws.Range("H5").Value = [double]"0.007588"
$ws.Range("M5").Value = [double]"37.74401233333333"
$ws.Range("N5").Value = [double]"113.232037"
$ws.Range("O5").Value = [double]"0.9981055287696164"
$ws.Range("P5").Value = [double]"0.9981055287696164"
$ws.Range("Q5").Value = [double]"0.09546718852844445"
$ws.Range("R5").Value = [double]"0.859204696756"
$ws.Range("S5").Value = [double]"0.9981055287696164"
$ws.Range("T5").Value = [double]"0.9981055287696164"
